$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -20.288
    6   = -22.186
    7   = -20.064
    8   = -21.972
    16  = -22.067
    20  = -20.375
    21  = -20.057
    28  = -21.918
    29  = -21.344
    30  = -21.915
    32  = -21.673
    40  = -20.003
    46  = -21.816
    51  = -21.95
    52  = -22.036
    57  = -22.208
    59  = -22.574
    62  = -22.098
    66  = -21.546
    73  = -20.082
    74  = -21.246
    77  = -20.603
    92  = -21.566
    100 = -22.165
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
